$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (Objetivos:) — data cells now hold the docente info instead of the
#     long "Possibilitar..." paragraph ---
$ws.Range("B10:C10").Value = "8188658 - Maria Auxiliadora Motta Barreto"

# --- Row 13 (was blank label / docente data) becomes "Programa resumido:" / "Semestral" ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13:C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14 (was "Programa resumido:") becomes "Short syllabus:" ---
$ws.Range("A14").Value = "Short syllabus:"
$shortSyllabus = @"
Introduction to Psychology applied to work.
Concept of Communication. 
Human Relations at Work. 
Psychology in Work Organizations .
Recruitment and Selection and 
Motivation
"@
$ws.Range("B14:C14").Value = $shortSyllabus

# --- Row 15 (was "Short syllabus:") becomes "Programa:" / "01/01/2012" ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15:C15").Value = "01/01/2012"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16 (was "Programa:") becomes "Syllabus:" ---
$ws.Range("A16").Value = "Syllabus:"
$syllabus = @"
Introduction : conceptualize psychology as science and application; psychology applied to work. The psychology of human relations at work.
Concept of Communication : Systems, functions , axioms of human communication. Communication processes and the social and communicative interaction in the company.
Human relations at work: the role of masks in human interaction; human relations in groups; how to be a part of a workgroup.
Psychology in Work Organizations : Organization concepts and work. Organization and work and their importance in mental health and worker productivity: stress, burnout , Karoshi syndrome ; L.E.R .; quality of life; sexual and moral harassment in the workplace ; alcohol and drugs at work; mental disorders in the company.
Recruitment and Selection: recruitment and selection of personnel ; placement and monitoring; performance evaluation; training and education; evaluation measures and their importance in the selection ; practical experiences in the classroom as facilitators in the selection process .
Motivation : the basic and psychological needs of human beings; motivation and conflicts; forgotten factors as motivators in the company : envy, jealousy , fear, abuse of power . Motivation Assessment.
"@
$ws.Range("B16:C16").Value = $syllabus

# --- Row 17 (was "Syllabus:") becomes "Avaliação:" and loses its B/C data ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").ClearContents()
$ws.Rows.Item(17).AutoFit()

# --- Row 18 (was "Avaliação:", blank data) becomes "Método:" / docente info ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18:C18").Value = "8188658 - Maria Auxiliadora Motta Barreto"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19 (was "Método:") becomes "Critério:" (data unchanged) ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20 (was "Critério:") becomes "Norma de recuperação:" (data unchanged) ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21 (was "Norma de recuperação:") becomes "Bibliografia:" (data unchanged) ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22 (was "Bibliografia:") becomes "Requisitos:" and loses its B/C data ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22:C22").ClearContents()
$ws.Rows.Item(22).AutoFit()

# --- Row 23 (was "Requisitos:", blank data) loses its A label, keeps the
#     requirement text moved up from row 24 ---
$ws.Range("A23").ClearContents()
$ws.Range("B23:C23").Value = "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# --- Row 24 no longer exists; delete it so the sheet ends at row 23 ---
$ws.Rows.Item(24).Delete()
